# Adapt column header formatting to respective input file names (#7)
#
# - "<name>_old"  -> "<name>_FV2310"
# - "<name>_new"  -> "<name>_FV2404"
# - Wrap the used range in an Excel Table ("Table1") whose column headers
#   mirror the renamed worksheet headers.
# - Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J => "<name>_FV2310"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($baseNames[$i])_FV2310"
}

# Column K ("diff") is left untouched.

# Columns L..U => "<name>_FV2404"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($baseNames[$i])_FV2404"
}

# Stash a copy of the header row's formatting out of the way (below the
# used range) so it survives the ClearFormats/ListObjects.Add round trip
# below without Excel baking a header-row dxf into the new table style.
$headerRange = $ws.Range("A1:U1")
$stashRange = $ws.Range("A59:U59")
$headerRange.Copy()
$stashRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

# Turn the used range into an Excel Table ("Table1") so the new header
# names are also reflected as table column headers.
$tableRange = $ws.Range("A1:U57")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"

# Restore the header row's original formatting.
$stashRange.Copy()
$headerRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$stashRange.Clear()

# Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
